# Built test rig, added new animation strcat's for use in test classes
#
# The sheet originally held a "fragment" of sample data that spanned
# A1:BN40. This trims it down to the smaller A1:BF27 rig used by the
# new test classes:
#   - row 1 only keeps its digits through column BF (BG1:BN1 cleared)
#   - the trailing filler rows 28-40 (AU/AV = 8 placeholders) are removed
#   - the active selection moves to K49

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim row 1 so it stops at column BF instead of BN.
$ws.Range("BG1:BN1").ClearContents()

# Remove the extra filler rows 28 through 40 entirely (shifts rows up).
$ws.Rows("28:40").Delete()

# Update the selected cell to match the new test rig layout.
$ws.Range("K49").Select()
